$d = $word.ActiveDocument

# The document ends with an empty paragraph (before the section break).
# Insert a brand-new paragraph right after it, then stamp it with the
# exact target OOXML (including the proofErr markers Word leaves around
# "The The") via InsertXML so the content matches the author's edit.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">The code is a set of unit tests for the Guess </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>The</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Number game. It checks if the game initializes correctly, if the guess method is working as expected when the player guesses a number higher or lower than the actual number, and if the game is won or lost when the player guesses the correct number or exceeds the maximum number of attempts.</w:t></w:r></w:p>'
$newPara.Range.InsertXML($newXml)
